# Fix some MRP-related issues
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the raw poll numbers that feed the formulas on the sheet.
$ws.Range("A2").Value = 32
$ws.Range("C2").Value = 13
$ws.Range("D2").Value = 8

$ws.Range("A10").Value = 46
$ws.Range("B10").Value = 48

# Reflect the last active selection left in the sheet after the edits.
$ws.Range("H13").Select()
